$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto price refresh diff.
# NumberFormat is forced to Text ("@") before assignment so that numeric-looking
# strings (e.g. "1.001", "0.07744") are preserved as literal text instead of being
# auto-converted to numbers by Excel, then the style is reset back to Normal so no
# stray formatting is introduced.
$updates = @(
    @{Row=2; Col=4; Value='30.580.52'},
    @{Row=2; Col=5; Value='  -1.53%  '},
    @{Row=3; Col=4; Value='1.925.85'},
    @{Row=3; Col=5; Value='  +1.06%  '},
    @{Row=4; Col=4; Value='1.001'},
    @{Row=4; Col=5; Value='  +0.18%  '},
    @{Row=5; Col=4; Value='246.37'},
    @{Row=5; Col=5; Value='  +2.20%  '},
    @{Row=6; Col=5; Value='  +0.16%  '},
    @{Row=7; Col=4; Value='0.4736'},
    @{Row=7; Col=5; Value='  -1.38%  '},
    @{Row=8; Col=4; Value='0.2927'},
    @{Row=8; Col=5; Value='  -1.94%  '},
    @{Row=9; Col=4; Value='0.06812'},
    @{Row=9; Col=5; Value='  +2.47%  '},
    @{Row=10; Col=4; Value='105.35'},
    @{Row=10; Col=5; Value='  +3.77%  '},
    @{Row=11; Col=4; Value='18.41'},
    @{Row=11; Col=5; Value='  -4.25%  '},
    @{Row=12; Col=2; Value='WrappedEther'},
    @{Row=12; Col=3; Value='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'},
    @{Row=12; Col=4; Value='1.919.57'},
    @{Row=12; Col=5; Value='  +0.74%  '},
    @{Row=13; Col=2; Value='TRON'},
    @{Row=13; Col=3; Value='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'},
    @{Row=13; Col=4; Value='0.07744'},
    @{Row=13; Col=5; Value='  +1.18%  '},
    @{Row=14; Col=4; Value='5.348'},
    @{Row=14; Col=5; Value='  +3.08%  '},
    @{Row=15; Col=4; Value='0.6717'},
    @{Row=15; Col=5; Value='  +0.55%  '},
    @{Row=16; Col=4; Value='286.88'},
    @{Row=16; Col=5; Value='  -6.88%  '},
    @{Row=17; Col=4; Value='30.623.37'},
    @{Row=17; Col=5; Value='  -1.33%  '},
    @{Row=18; Col=4; Value='13.07'},
    @{Row=18; Col=5; Value='  -1.61%  '},
    @{Row=19; Col=4; Value='0.000007655'},
    @{Row=19; Col=5; Value='  +0.70%  '},
    @{Row=20; Col=4; Value='1.000'},
    @{Row=20; Col=5; Value='  +0.06%  '},
    @{Row=21; Col=4; Value='2.164.84'},
    @{Row=21; Col=5; Value='  +1.24%  '},
    @{Row=22; Col=4; Value='5.428'},
    @{Row=22; Col=5; Value='  +3.59%  '},
    @{Row=23; Col=5; Value='  +0.15%  '},
    @{Row=24; Col=4; Value='6.270'},
    @{Row=24; Col=5; Value='  +0.12%  '},
    @{Row=25; Col=4; Value='9.409'},
    @{Row=25; Col=5; Value='  +0.04%  '},
    @{Row=26; Col=4; Value='168.62'},
    @{Row=26; Col=5; Value='  +0.04%  '},
    @{Row=27; Col=4; Value='20.73'},
    @{Row=27; Col=5; Value='  -0.76%  '},
    @{Row=28; Col=4; Value='2.137'},
    @{Row=28; Col=5; Value='  +6.28%  '},
    @{Row=29; Col=4; Value='0.1087'},
    @{Row=29; Col=5; Value='  -2.72%  '},
    @{Row=30; Col=5; Value='  +0.27%  '},
    @{Row=31; Col=4; Value='4.169'},
    @{Row=31; Col=5; Value='  -0.67%  '},
    @{Row=32; Col=4; Value='4.007'},
    @{Row=32; Col=5; Value='  -0.60%  '},
    @{Row=33; Col=4; Value='0.05070'},
    @{Row=33; Col=5; Value='  -0.80%  '},
    @{Row=34; Col=4; Value='0.7399'},
    @{Row=34; Col=5; Value='  -1.94%  '},
    @{Row=35; Col=4; Value='1.153'},
    @{Row=35; Col=5; Value='  -1.94%  '},
    @{Row=36; Col=4; Value='0.02088'},
    @{Row=36; Col=5; Value='  +3.50%  '},
    @{Row=37; Col=5; Value='  -1.00%  '},
    @{Row=38; Col=4; Value='2.694'},
    @{Row=38; Col=5; Value='  -1.11%  '},
    @{Row=39; Col=4; Value='2.063'},
    @{Row=39; Col=5; Value='  -0.81%  '},
    @{Row=40; Col=4; Value='111.05'},
    @{Row=40; Col=5; Value='  +1.30%  '},
    @{Row=41; Col=4; Value='0.8765'},
    @{Row=41; Col=5; Value='  -1.22%  '},
    @{Row=42; Col=4; Value='0.4439'},
    @{Row=42; Col=5; Value='  +4.31%  '},
    @{Row=43; Col=4; Value='5.922'},
    @{Row=43; Col=5; Value='  +3.00%  '},
    @{Row=44; Col=5; Value='  +0.18%  '},
    @{Row=45; Col=4; Value='67.63'},
    @{Row=45; Col=5; Value='  -1.17%  '},
    @{Row=46; Col=4; Value='7.269'},
    @{Row=46; Col=5; Value='  -2.00%  '},
    @{Row=47; Col=4; Value='9.342'},
    @{Row=47; Col=5; Value='  +0.88%  '},
    @{Row=48; Col=2; Value='Algorand'},
    @{Row=48; Col=3; Value='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'},
    @{Row=48; Col=4; Value='0.1240'},
    @{Row=48; Col=5; Value='  -0.09%  '},
    @{Row=49; Col=2; Value='Decentraland'},
    @{Row=49; Col=3; Value='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'},
    @{Row=49; Col=4; Value='0.4124'},
    @{Row=49; Col=5; Value='  +6.06%  '},
    @{Row=50; Col=2; Value='Elrond'},
    @{Row=50; Col=3; Value='https://coinranking.com/coin/omwkOTglq+elrond-egld'},
    @{Row=50; Col=4; Value='35.25'},
    @{Row=50; Col=5; Value='  +0.47%  '},
    @{Row=51; Col=2; Value='BitcoinSV'},
    @{Row=51; Col=3; Value='https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'},
    @{Row=51; Col=4; Value='47.03'},
    @{Row=51; Col=5; Value='  +8.99%  '}
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

